{"js": "// Apply the benchmark-table value corrections described by the diff.\n// The document is a single table, one column, one row per metric value.\n// Rows are addressed 0-based via table.rows.items[i].cells.items[0].\nconst table = context.document.body.tables.getFirst();\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Map of 0-based row index -> new cell text (replaces ALL runs/tabs in\n// that cell's paragraph with a single run carrying the new text).\nconst updates = {\n  0: \"0M\",       // was \"100\"\n  1: \"0M\",       // was \"0\"\n  2: \"0M\",       // was \"2397\"\n  3: \"78\",       // was \"3\"\n  4: \"0.00003\",  // was \"0.00004\"\n  5: \"0.00022\",  // was \"0.00005\"\n  11: \"0.00371\", // was \"0.00013\"\n  43: \"100\",     // was the \"1 <tab> 0.00003 ... 100.0\" run sequence\n  44: \"0\",       // was the \"64 <tab> 0.00003 ... 100.0\" run sequence\n  45: \"2397\",    // was the \"10 <tab> 0.00003 ... 100.0\" run sequence\n};\n\nfor (const [idx, text] of Object.entries(updates)) {\n  const cell = rows.items[Number(idx)].cells.items[0];\n  cell.value = text;\n}\n\nawait context.sync();\n", "ps1": "# Apply the benchmark-table value corrections described by the diff.\n# The document is a single table, one column, one row per metric value.\n# Word COM collections (Rows/Cells) are 1-based.\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n$t.Cell(1, 1).Range.Text = \"0M\"        # was \"100\"\n$t.Cell(2, 1).Range.Text = \"0M\"        # was \"0\"\n$t.Cell(3, 1).Range.Text = \"0M\"        # was \"2397\"\n$t.Cell(4, 1).Range.Text = \"78\"        # was \"3\"\n$t.Cell(5, 1).Range.Text = \"0.00003\"   # was \"0.00004\"\n$t.Cell(6, 1).Range.Text = \"0.00022\"   # was \"0.00005\"\n$t.Cell(12, 1).Range.Text = \"0.00371\"  # was \"0.00013\"\n$t.Cell(44, 1).Range.Text = \"100\"      # was \"1<tab>0.00003...100.0\" run sequence\n$t.Cell(45, 1).Range.Text = \"0\"        # was \"64<tab>0.00003...100.0\" run sequence\n$t.Cell(46, 1).Range.Text = \"2397\"     # was \"10<tab>0.00003...100.0\" run sequence\n"}
